# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet, shifting the existing N/O/P columns
# (Late, heading/Date, Outstanding) one place to the right, and
# switch the active sheet/selection from "Summary" to
# "Repayment schedule".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N; existing N->O, O->P, P->Q.
$ws.Columns("N:N").Insert()

# The new column inherits the (approx.) width of the column that used
# to sit there before the shift ("In Advance", column M).
$ws.Columns("N:N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet (this also clears
# tabSelected on whichever sheet was previously active, e.g. "Summary").
$ws.Activate()

# Update the selection on the newly active sheet.
$ws.Range("K19").Select()
